$wb = $excel.ActiveWorkbook

# The "createRegionalEvent" sheet holds event names; append two new test
# automation event name values below the existing data (A1:B2 -> A1:B4).
$ws = $wb.Worksheets.Item("createRegionalEvent")

$ws.Range("A3").Value = "TestAutomation20191217141856"
$ws.Range("A4").Value = "TestAutomation20191217142208"
